$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    (-0.2095100884369551,-0.8379002189740954),
    (-0.4671905040515311,-1.230702036034245),
    (-0.5455376767895779,-1.607184072506625),
    (-0.2032515484787179,-0.6196333778449691),
    (-0.2015325903169599,-0.2609675993684708),
    (-0.2371420086554795,-0.8937657444083623),
    (-0.4809405008460905,-1.104493371240874),
    (-0.6462027447962853,-1.722176332716887),
    (-0.4984316114019968,-0.4374574618321648),
    (-0.2698359325064053,-0.9460183002241354),
    (-0.2200911257904904,-0.4229785037097781),
    (-0.1037873042165333,-0.2145148933534846),
    (-0.2724062589405265,-0.3162005862098188),
    (-0.06732061235437686,0.1255420728274558),
    (-0.2396914543690218,-0.6283176742554456),
    (0.003036054274278797,0.4451127093981985),
    (0.02156687078558003,0.5551525444481236),
    (0.05624229038695133,0.3490889584933873),
    (-0.158246147902445,-0.4488377904159373),
    (-0.03253682032423353,0.2490004173988481),
    (-0.07423985150057538,0.2856838034772378),
    (0.02694808364132478,0.4948601289525356),
    (0.8451072203291963,2.384501066679122),
    (0.1288135216087946,0.6573290021976542),
    (0.1839054168366249,0.6116695573938347),
    (0.1306071906351134,0.6915887987496229),
    (0.3723775497480384,1.27457830695872),
    (0.7021712938188642,2.299311468773102),
    (0.2511787273034483,0.7904614385661362),
    (0.08188460933219721,0.5947565271340425),
    (0.2328830996605502,1.036557800243321),
    (0.1516001679706656,0.7873743958802361),
    (0.1085785059715446,0.4817198996339583),
    (0.4059551404095935,1.183761674201925),
    (0.2855661592284053,0.7382262801659634),
    (-0.02985859994145716,-0.04644764224621437),
    (0.3606759753309736,1.655840245639183),
    (-0.1134971026603982,-0.2845369871992774),
    (0.1485976969633974,0.6291346820648994),
    (-0.1666981911741042,0.5775014000248526),
    (0.2832620069664291,1.57331484963083),
    (0.113888612116753,0.8138681769111067),
    (0.04251601851196032,0.3601140775006053),
    (-0.1212609069821894,0.03488624435914889),
    (-0.1721606371902295,-0.4408384506887962),
    (-0.1676381673409932,-0.350346897799793),
    (-0.2211252889132141,-0.5384615165485332),
    (-0.2238407406011904,-0.5948442133249261),
    (-0.1008528053657672,-0.1488916354239387),
    (-0.2284903471496285,-0.630849053808473),
    (-0.2284903471496285,-0.630849053808473),
    (-0.1953346499633997,-0.3891390135927847),
    (-0.2183101419596003,-0.5423375210844638),
    (-0.1549616035746775,-0.3375691409173693),
    (-0.1163979145544606,-0.2837291266230884),
    (-0.1832469983269261,-0.3404896288290168),
    (-0.145545391540654,-0.2601452411715781),
    (-0.2223132843059765,-0.6492027115795734),
    (-0.2536643487663319,-0.7110722138153006),
    (-0.2683979652135061,-0.3305365139864404),
    (-0.1491951328431063,-0.07310099463401706),
    (-0.4748670785773039,-1.283841480614283),
    (-0.2972107521500132,-0.6572935007319717),
    (-0.2483761848036102,-0.6116519566086053),
    (-0.1037779011186861,-0.03484158900219716),
    (0.02302464629094756,0.1812275194303561)
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $pair = $values[$i]
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
}
